# LOQ4085.xlsx edit
#
# The upstream data feeding this sheet shifted by one row: the row that used
# to hold only the "8151869 - Livia Chaguri e Carvalho" value (row 13, with no
# label in column A) disappears, and several of the remaining label rows
# (column A keeps its caption) now show the value that used to belong to a
# different row. We reproduce this by first copying every affected value into
# its new home (using values-only paste so existing cell formatting/style is
# left untouched), then deleting the now-redundant row, and finally filling in
# the one genuinely new piece of text ("Semestral").
#
# All row/column references below use the ORIGINAL (pre-edit) row numbers,
# since the copies are all performed before the row delete happens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# IMPORTANT: several of these moves form a chain (row 19 -> row 20 -> row 21
# -> row 22, with row 13 also feeding rows 10/19, and row 8 feeding row 16).
# Each source row must be read out *before* it is itself overwritten as the
# destination of an earlier move, so apply the moves starting from the end of
# the chain and working backwards.

# Row 22 ("Bibliografia:" -- becomes row 21) now shows the recuperação-average
# text that used to live in row 21.
$ws.Range("B21").Copy() | Out-Null
$ws.Range("B22").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C22").PasteSpecial($xlPasteValues) | Out-Null

# Row 21 ("Norma de recuperação:" -- becomes row 20) now shows the grading
# criteria text that used to live in row 20.
$ws.Range("B20").Copy() | Out-Null
$ws.Range("B21").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("C20").Copy() | Out-Null
$ws.Range("C21").PasteSpecial($xlPasteValues) | Out-Null

# Row 20 ("Critério:" -- becomes row 19) now shows the evaluation-method text
# that used to live in row 19.
$ws.Range("B19").Copy() | Out-Null
$ws.Range("B20").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("C19").Copy() | Out-Null
$ws.Range("C20").PasteSpecial($xlPasteValues) | Out-Null

# Row 19 ("Método:" -- becomes row 18 once row 13 is removed) picks up the
# docente value from row 13.
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B19").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C19").PasteSpecial($xlPasteValues) | Out-Null

# Row 10 ("Objetivos:") picks up the docente value that used to live in row 13.
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B10").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C10").PasteSpecial($xlPasteValues) | Out-Null

# Row 16 ("Programa:" -- becomes row 15) now shows the activation date that
# used to live in row 8.
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B16").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C16").PasteSpecial($xlPasteValues) | Out-Null

$excel.CutCopyMode = $false

# Remove the now-redundant standalone value row 13 (B/C = docente value, no
# column A label). Everything below shifts up by one row.
$ws.Rows(13).EntireRow.Delete()

# Row 13 ("Programa resumido:", after the shift) gets the genuinely new value
# "Semestral", with the same 60pt row height used by its neighbors.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# Row heights for rows 15 and 21 change from 60pt to 120pt in the new layout.
$ws.Rows(15).RowHeight = 120
$ws.Rows(21).RowHeight = 120
